$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 0.0009722222222222222
$ws.Range("K2").Value = 4157
$ws.Range("L2").Value = 0.008314
